# Apply the "Code VS" URL update + "Value Types" simplification for the
# observations-summary sheet (fhir/ig/tei 0.2.1).
#
# Summary of the change:
#  - The "http://hl7.org/fhir/ValueSet/observation-codes (example)" Code VS
#    used by several rows is replaced by a new, more specific ValueSet:
#    "https://interoperabilidad.minsal.cl/fhir/ig/tei/ValueSet/TipoDeObservacion (extensible)"
#  - The "null#108217004" code on the Anamnesis row becomes "null#84100007"
#  - The Indice Comorbilidad row's Code VS (previously VSIndicecomorbilidad)
#    now also points to the TipoDeObservacion ValueSet, and its Value Types
#    cell collapses from a long list down to just "CodeableConceptĵ"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newCodeVs = "https://interoperabilidad.minsal.cl/fhir/ig/tei/ValueSet/TipoDeObservacion (extensible)"

# Row 2 - ObservationAnamnesisLE
$ws.Range("E2").Value = "null#84100007"
$ws.Range("F2").Value = $newCodeVs

# Row 3 - ObservationDiscapacidadLE
$ws.Range("F3").Value = $newCodeVs

# Row 4 - ObservationIndiceComorbilidadLE
$ws.Range("F4").Value = $newCodeVs
$ws.Range("H4").Value = "CodeableConceptĵ"

# Row 5 - ObservationIniciarCuidadorLE
$ws.Range("F5").Value = $newCodeVs
